$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (capex functionality recalculated costs)
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 139.047619047619
$ws.Range("C3").Value = 2380.952380952381

# Add new row for biomass fuel type
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "biomass"
$ws.Range("C4").Value = 1000
